# The data table (colo list) gets a new data center entry "AKX" (Aktobe,
# Kazakhstan) inserted as a new row right before the existing "IAD" row
# (originally row 276), pushing all following rows down by one and growing
# the used range from A1:H336 to A1:H337.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 276; everything from the old row 276
# onward (IAD, ATL, ... YHZ) shifts down to 277..337.
$ws.Rows.Item(276).Insert()

# The freshly inserted row has no formatting of its own; copy the
# column-A style (bold text + thin border) from the row right below it
# (which is the original "IAD" row, now at 277) so the new "AKX" cell
# matches the look of every other colo-code cell in column A.
$ws.Cells.Item(277, 1).Copy()
$ws.Cells.Item(276, 1).PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the Aktobe, Kazakhstan colo record.
$ws.Cells.Item(276, 1).Value = "AKX"
$ws.Cells.Item(276, 2).Value = "Aktobe, Kazakhstan"
$ws.Cells.Item(276, 3).Value = "Asia Pacific"
$ws.Cells.Item(276, 4).Value = "Aktobe"
$ws.Cells.Item(276, 5).Value = "Kazakhstan"
$ws.Cells.Item(276, 6).Value = "KZ"
$ws.Cells.Item(276, 7).Value = 50.286922
$ws.Cells.Item(276, 8).Value = 57.224121
